$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 657.7217000000001
$ws.Range("I15").Value = 657.7217000000001
$ws.Range("K15").Value = 1973.1651
$ws.Range("M15").Value = -1804.1651

$ws.Range("H137").Value = 5000846.5
$ws.Range("I137").Value = 638.76
$ws.Range("J137").Value = 13334526
$ws.Range("K137").Value = 1916.28
$ws.Range("L137").Value = 40003578
$ws.Range("M137").Value = 633.72
$ws.Range("N137").Value = -40008678

$ws.Range("H138").Value = 4116631.5
$ws.Range("I138").Value = 6290016
$ws.Range("J138").Value = 2725
$ws.Range("K138").Value = 18870048
$ws.Range("L138").Value = 8175
$ws.Range("M138").Value = -18864908
$ws.Range("N138").Value = -18455

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 195000
$ws.Range("J140").Value = 195000
$ws.Range("L140").Value = 195000
$ws.Range("N140").Value = -205360

$ws.Range("H141").Value = 888.86536
$ws.Range("I141").Value = 888.86536
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2666.59608
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2513.40392
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5210636
$ws.Range("I132").Value = 9261485
$ws.Range("J132").Value = 2401.9048
$ws.Range("K132").Value = 27784455
$ws.Range("L132").Value = 7205.714399999999
$ws.Range("M132").Value = -27781925
$ws.Range("N132").Value = -12265.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12353431
$ws.Range("I31").Value = 10553.467
$ws.Range("K31").Value = 10553.467
$ws.Range("M31").Value = -10258.467

$ws.Range("H34").Value = 12353431
$ws.Range("I34").Value = 10553.467
$ws.Range("K34").Value = 10553.467
$ws.Range("M34").Value = -10351.467

$ws.Range("H58").Value = 1132.9778
$ws.Range("I58").Value = 523.54285
$ws.Range("J58").Value = 3266
$ws.Range("K58").Value = 523.54285
$ws.Range("L58").Value = 3266
$ws.Range("M58").Value = -320.54285
$ws.Range("N58").Value = -3672

$ws.Range("H136").Value = 1132.9778
$ws.Range("I136").Value = 523.54285
$ws.Range("J136").Value = 3266
$ws.Range("K136").Value = 1570.62855
$ws.Range("L136").Value = 9798
$ws.Range("M136").Value = 979.3714499999999
$ws.Range("N136").Value = -14898

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 809.37
$ws.Range("I131").Value = 471.92307
$ws.Range("J131").Value = 859.7931
$ws.Range("K131").Value = 1415.76921
$ws.Range("L131").Value = 2579.3793
$ws.Range("M131").Value = 3624.23079
$ws.Range("N131").Value = -12659.3793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3033.0876
$ws.Range("I132").Value = 2205.9023
$ws.Range("J132").Value = 5152.75
$ws.Range("K132").Value = 6617.706900000001
$ws.Range("L132").Value = 15458.25
$ws.Range("M132").Value = -4087.706900000001
$ws.Range("N132").Value = -20518.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2235.1875
$ws.Range("I16").Value = 2069.182
$ws.Range("J16").Value = 2600.4
$ws.Range("K16").Value = 2069.182
$ws.Range("L16").Value = 2600.4
$ws.Range("M16").Value = -1899.182
$ws.Range("N16").Value = -2940.4

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 7048645.5
$ws.Range("I132").Value = 3939.8867
$ws.Range("K132").Value = 11819.6601
$ws.Range("M132").Value = -9289.660100000001

$ws.Range("H136").Value = 8336123
$ws.Range("I136").Value = 9435000
$ws.Range("J136").Value = 16057.857
$ws.Range("K136").Value = 28305000
$ws.Range("L136").Value = 48173.571
$ws.Range("M136").Value = -28302450
$ws.Range("N136").Value = -53273.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2111.1924
$ws.Range("I122").Value = 2163
$ws.Range("J122").Value = 1826.25
$ws.Range("K122").Value = 6489
$ws.Range("L122").Value = 5478.75
$ws.Range("M122").Value = -4039
$ws.Range("N122").Value = -10378.75

$ws.Range("H125").Value = 60715
$ws.Range("J125").Value = 60715
$ws.Range("L125").Value = 60715
$ws.Range("N125").Value = -70555

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040

$ws.Range("H131").Value = 75999
$ws.Range("J131").Value = 75999
$ws.Range("L131").Value = 75999
$ws.Range("N131").Value = -86079

$ws.Range("H132").Value = 1846.9246
$ws.Range("I132").Value = 1610.0264
$ws.Range("J132").Value = 2447.0667
$ws.Range("K132").Value = 4830.0792
$ws.Range("L132").Value = 7341.2001
$ws.Range("M132").Value = -2300.0792
$ws.Range("N132").Value = -12401.2001

$ws.Range("H136").Value = 1320.4546
$ws.Range("I136").Value = 1094.84
$ws.Range("J136").Value = 2025.5
$ws.Range("K136").Value = 3284.52
$ws.Range("L136").Value = 6076.5
$ws.Range("M136").Value = -734.5199999999995
$ws.Range("N136").Value = -11176.5

$ws.Range("H137").Value = 61143.332
$ws.Range("J137").Value = 61143.332
$ws.Range("L137").Value = 61143.332
$ws.Range("N137").Value = -71343.33199999999

$ws.Range("H140").Value = 75167.5
$ws.Range("J140").Value = 75167.5
$ws.Range("L140").Value = 75167.5
$ws.Range("N140").Value = -85527.5
